# Update "想去人数" (interest count) figures on the 展览, 演出, and 全部类型
# sheets to their refreshed values (output regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 224
$ws.Range("F3").Value = 373
$ws.Range("F5").Value = 22
$ws.Range("F6").Value = 81
$ws.Range("F8").Value = 379
$ws.Range("F9").Value = 4587
$ws.Range("F10").Value = 4587
$ws.Range("F12").Value = 444
$ws.Range("F13").Value = 1079
$ws.Range("F14").Value = 595
$ws.Range("F15").Value = 4080
$ws.Range("F16").Value = 155
$ws.Range("F17").Value = 154
$ws.Range("F18").Value = 42
$ws.Range("F19").Value = 196
$ws.Range("F20").Value = 3379
$ws.Range("F24").Value = 2916
$ws.Range("F25").Value = 118
$ws.Range("F27").Value = 6
$ws.Range("F28").Value = 138
$ws.Range("F29").Value = 168
$ws.Range("F30").Value = 169
$ws.Range("F32").Value = 45
$ws.Range("F36").Value = 5283
$ws.Range("F37").Value = 721
$ws.Range("F38").Value = 380
$ws.Range("F39").Value = 77
$ws.Range("F41").Value = 12
$ws.Range("F42").Value = 1039
$ws.Range("F43").Value = 418
$ws.Range("F44").Value = 19
$ws.Range("F45").Value = 1921
$ws.Range("F47").Value = 51
$ws.Range("F48").Value = 684
$ws.Range("F49").Value = 816

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F21").Value = 712

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 224
$ws.Range("F6").Value = 22
$ws.Range("F9").Value = 379
$ws.Range("F10").Value = 4587
$ws.Range("F11").Value = 4587
$ws.Range("F16").Value = 444
$ws.Range("F17").Value = 1079
$ws.Range("F18").Value = 595
$ws.Range("F19").Value = 4080
$ws.Range("F20").Value = 155
$ws.Range("F21").Value = 154
$ws.Range("F22").Value = 196
$ws.Range("F23").Value = 3379
$ws.Range("F24").Value = 2916
$ws.Range("F25").Value = 118
$ws.Range("F27").Value = 138
$ws.Range("F28").Value = 168
$ws.Range("F29").Value = 169
$ws.Range("F36").Value = 5283
$ws.Range("F38").Value = 721
$ws.Range("F39").Value = 380
$ws.Range("F41").Value = 77
$ws.Range("F43").Value = 1039
$ws.Range("F44").Value = 418
$ws.Range("F45").Value = 19
$ws.Range("F46").Value = 1921
$ws.Range("F48").Value = 51
$ws.Range("F49").Value = 684
$ws.Range("F50").Value = 816
